$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "59.132.80"
$ws.Range("E2").Value = "  +0.31%  "

$ws.Range("D3").Value = "2.518.34"
$ws.Range("E3").Value = "  +0.69%  "

$ws.Range("E4").Value = "  +0.09%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "536.43"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.07%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "139.80"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -2.52%  "

$ws.Range("E7").Value = "  +0.37%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.561"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -1.55%  "

$ws.Range("D9").Value = "2.520.07"
$ws.Range("E9").Value = "  -0.29%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0995"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.12%  "

$ws.Range("E11").Value = "  +1.46%  "

$ws.Range("E12").Value = "  -0.98%  "

$ws.Range("E13").Value = "  +1.60%  "

$ws.Range("D14").Value = "2.964.89"
$ws.Range("E14").Value = "  +0.94%  "

$ws.Range("D15").Value = "59.101.91"
$ws.Range("E15").Value = "  +0.54%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "22.91"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -2.26%  "

$ws.Range("E17").Value = "  +1.21%  "

$ws.Range("D18").Value = "2.522.86"
$ws.Range("E18").Value = "  +0.33%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "10.89"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -2.48%  "

$ws.Range("E20").Value = "  -0.81%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "321.45"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.18%  "

$ws.Range("E22").Value = "  -0.27%  "

$ws.Range("E23").Value = "  +1.44%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "62.42"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.17%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.422"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -2.80%  "

$ws.Range("E26").Value = "  +1.42%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.999"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.59%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.75"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.56%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "6.72"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.78%  "

$ws.Range("E30").Value = "  +0.46%  "

$ws.Range("D31").Value = "0.0₃0765"
$ws.Range("E31").Value = "  -0.07%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "160.75"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.81%  "

$ws.Range("E33").Value = "  +0.29%  "

$ws.Range("E34").Value = "  -4.56%  "

$ws.Range("E35").Value = "  +2.32%  "

$ws.Range("E36").Value = "  -0.45%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "4.19"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -2.84%  "

$ws.Range("E38").Value = "  -1.85%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "36.94"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +1.60%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.63"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.16%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.802"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.63%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "5.25"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -4.70%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "282.54"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -4.96%  "

$ws.Range("E44").Value = "  +0.36%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "10.88"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +1.04%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.594"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -1.44%  "

$ws.Range("E47").Value = "  +0.21%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "122.28"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -1.86%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "18.48"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.67%  "

$ws.Range("E50").Value = "  -0.10%  "

$ws.Range("E51").Value = "  -1.88%  "
